# Rollover the collection data: swap out the three "past" collections
# (Print Paradise, Seeing Stars, Vitamin Sea) for the new ones being
# rolled in (Meet Your Match, Tough Luxe, Power Surge) along with their
# new product ids, while leaving the new_title/new_price/new_template
# columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Meet Your Match 3 Items"
$ws.Range("B2").Value = 1371333656668

$ws.Range("A3").Value = "Tough Luxe 3 Items"
$ws.Range("B3").Value = 1371334639708

$ws.Range("A4").Value = "Power Surge 3 Items"
$ws.Range("B4").Value = 1372047638620
